# Generate Report for Handoff
# Replace the old run's file-id / content-hash with the new run's, refresh
# the handoff/handback timestamps, and blank out the (not-yet-produced)
# handback info for both locales.

$wb = $excel.ActiveWorkbook

$oldId   = "0aa75caa-44d8-456c-9b37-6204ba854d46"
$newId   = "6013bfae-8619-4a51-9d7c-53a3107a007d"
$newHash = "4c7641e7709e872872be76be44726113e5ced95c"

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Overview sheet: file name / path columns + the HO xliff generation time
# ---------------------------------------------------------------------
$ws1.Range("A2").Value = ($newId + ".md")
$ws1.Range("B2").Value = ("e2e\" + $newId + ".md")
$ws1.Range("G2").Value = "2016-08-16 22:58:33"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0b5cbad21c79a88834584ad4f7c8b570845cef92/e2e/" + $newId + ".md", [System.Type]::Missing, [System.Type]::Missing, ("e2e\" + $newId + ".md"))

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws2.Range("A2").Value = ($newId + ".md")
$ws2.Range("G2").Value = ($newId + "." + $newHash + ".zh-cn.xlf")
$ws2.Range("H2").Value = "2016-08-16 22:58:28"
$ws2.Range("I2").Value = ""
$ws2.Range("I2").Style = "Normal"
$ws2.Range("J2").Value = ""
$ws2.Range("K2").Value = "0001-01-01 00:00:00"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0b5cbad21c79a88834584ad4f7c8b570845cef92/e2e/" + $newId + ".md", [System.Type]::Missing, [System.Type]::Missing, ($newId + ".md"))

$ws2.Columns.Item(9).ColumnWidth = 17.8333333333333
$ws2.Columns.Item(10).ColumnWidth = 20.8333333333333

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws3.Range("A2").Value = ($newId + ".md")
$ws3.Range("G2").Value = ($newId + "." + $newHash + ".de-de.xlf")
$ws3.Range("H2").Value = "2016-08-16 22:58:33"
$ws3.Range("I2").Value = ""
$ws3.Range("I2").Style = "Normal"
$ws3.Range("J2").Value = ""
$ws3.Range("K2").Value = "0001-01-01 00:00:00"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0b5cbad21c79a88834584ad4f7c8b570845cef92/e2e/" + $newId + ".md", [System.Type]::Missing, [System.Type]::Missing, ($newId + ".md"))

$ws3.Columns.Item(9).ColumnWidth = 17.8333333333333
$ws3.Columns.Item(10).ColumnWidth = 20.8333333333333

Write-Output "Report updated for handoff."
